$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new payment entry: "Check 7/7/16" = 300.18 (rows 13-14, column H) ---
# Added first so its shared string index precedes the "Café Espresso" one,
# matching the order new strings appear in the target sharedStrings table.
$ws.Range("H13").Value = "Check 7/7/16"
$ws.Range("H14").Value = 300.18

# --- Add new purchase category "Café Espresso" = 19.71 (rows 7-8, column E) ---
$ws.Range("E7").Value = "Café Espresso"
$ws.Range("E8").Value = 19.71

# --- Update percentage-split formulas in row 2 ---
# D2 share drops from 35% to 30%, F2 share rises from 10% to 15%
$ws.Range("D2").Formula = "=H2 * 0.3"
$ws.Range("F2").Formula = "=H2*0.15"

# --- Re-enter the row 4 and row 5 totals as one fill across A:F so the ---
# --- engine groups them into shared formulas, matching the target layout ---
$ws.Range("A4:F4").Formula = "=SUMPRODUCT(A7:A53,MOD(ROW(A7:A53)+1,2))"
$ws.Range("A5:F5").Formula = "=A2-A4"

# --- Move the active selection to K14 ---
$ws.Range("K14").Select()
